# CIERRE 15 FEB 22
# Fill in the credit/payment rows for February 2022 (rows 5-15) on the
# "REMISIONES FEBRERO 2022" sheet, extend the running-folio formula in
# column B down through row 36, and move the active selection to H12.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("REMISIONES FEBRERO 2022")

# --- Rows 5-15: date recibida (A), cliente (D), importe (E), fecha de pago (F), importe pagado (G) ---

$ws.Range("A5").Value = 44596
$ws.Range("D5").Value = "EL PRIMO"
$ws.Range("E5").Value = 15980
$ws.Range("F5").Value = 44598
$ws.Range("G5").Value = 15980

$ws.Range("A6").Value = 44597
$ws.Range("D6").Value = "OBRADOR"
$ws.Range("E6").Value = 288
$ws.Range("F6").Value = 44597
$ws.Range("G6").Value = 288

$ws.Range("A7").Value = 44597
$ws.Range("D7").Value = "OBRADOR"
$ws.Range("E7").Value = 500
$ws.Range("F7").Value = 44597
$ws.Range("G7").Value = 500

$ws.Range("A8").Value = 44600
$ws.Range("D8").Value = "OBRADOR"
$ws.Range("E8").Value = 263
$ws.Range("F8").Value = 44601
$ws.Range("G8").Value = 263

$ws.Range("A9").Value = 44600
$ws.Range("D9").Value = "HERRADURA GUSTAVO"
$ws.Range("E9").Value = 3726
$ws.Range("F9").Value = 44600
$ws.Range("G9").Value = 3726

$ws.Range("A10").Value = 44600
$ws.Range("D10").Value = "GABRIEL"
$ws.Range("E10").Value = 2704
$ws.Range("F10").Value = 44601
$ws.Range("G10").Value = 2704

$ws.Range("A11").Value = 44601
$ws.Range("D11").Value = "OBRADOR"
$ws.Range("E11").Value = 184
$ws.Range("F11").Value = 44601
$ws.Range("G11").Value = 184

$ws.Range("A12").Value = 44603
$ws.Range("D12").Value = "EL PRIMO"
$ws.Range("E12").Value = 18178
$ws.Range("F12").Value = 44605
$ws.Range("G12").Value = 18178

$ws.Range("A13").Value = 44603
$ws.Range("D13").Value = "OBRADOR"
$ws.Range("E13").Value = 78
$ws.Range("F13").Value = 44604
$ws.Range("G13").Value = 78

$ws.Range("A14").Value = 44604
$ws.Range("D14").Value = "OBRADOR"
$ws.Range("E14").Value = 212
$ws.Range("F14").Value = 44604
$ws.Range("G14").Value = 212

# Row 15 is still unpaid: no Fecha de pago / Importe D/Pago yet, so F15/G15
# stay blank and H15 (=E15-G15) carries the full 506 as outstanding balance.
$ws.Range("A15").Value = 44606
$ws.Range("D15").Value = "OBRADOR"
$ws.Range("E15").Value = 506

# --- Extend the running folio number (column B) from row 18 down to row 36 ---

for ($r = 18; $r -le 36; $r++) {
    $prev = $r - 1
    $ws.Range("B$r").Formula = "=B$prev+1"
}

# --- Restore the active-cell selection as left by the editor ---

$ws.Activate()
$ws.Range("H12").Select()
